$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2023-04-27, serial 45043) needs to be
# inserted as the second data row (row 10), pushing every existing data
# row (old rows 10-125) down by one (new rows 11-126).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new observation.
$ws.Cells.Item(10, 1).Value = 6
$ws.Cells.Item(10, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 45043
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 220
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 19000
$ws.Cells.Item(10, 13).Value = 17909
$ws.Cells.Item(10, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 1378
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"
